$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the manager_id value for row 1 (was the text "Null", causing a type
# mismatch with the numeric manager_id values used elsewhere in the column).
$ws.Range("F1").ClearContents()
